$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44169
$ws.Cells.Item(2, 8).Value = "Verde"
$ws.Cells.Item(2, 10).Value = 600
$ws.Cells.Item(2, 11).Value = 1600
$ws.Cells.Item(2, 12).Value = 1600
$ws.Cells.Item(2, 13).Value = 1600
$ws.Cells.Item(2, 16).Value = 1600

$ws.Cells.Item(3, 4).Value = 44525
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 10).Value = 180
$ws.Cells.Item(3, 11).Value = 1600
$ws.Cells.Item(3, 12).Value = 1600
$ws.Cells.Item(3, 13).Value = 1600
$ws.Cells.Item(3, 16).Value = 1600

$ws.Cells.Item(4, 4).Value = 44518
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 10).Value = 180
$ws.Cells.Item(4, 11).Value = 1600
$ws.Cells.Item(4, 12).Value = 1600
$ws.Cells.Item(4, 13).Value = 1600
$ws.Cells.Item(4, 16).Value = 1600

$ws.Cells.Item(5, 4).Value = 44487
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 10).Value = 120
$ws.Cells.Item(5, 11).Value = 1800
$ws.Cells.Item(5, 12).Value = 1800
$ws.Cells.Item(5, 13).Value = 1800
$ws.Cells.Item(5, 16).Value = 1800

$ws.Cells.Item(6, 4).Value = 44474
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 10).Value = 250
$ws.Cells.Item(6, 11).Value = 2000
$ws.Cells.Item(6, 12).Value = 2000
$ws.Cells.Item(6, 13).Value = 2000
$ws.Cells.Item(6, 16).Value = 2000

$ws.Cells.Item(7, 4).Value = 44532
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 10).Value = 180
$ws.Cells.Item(7, 11).Value = 1500
$ws.Cells.Item(7, 12).Value = 1500
$ws.Cells.Item(7, 13).Value = 1500
$ws.Cells.Item(7, 16).Value = 1500

$ws.Cells.Item(8, 4).Value = 44497
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 10).Value = 120
$ws.Cells.Item(8, 11).Value = 1800
$ws.Cells.Item(8, 12).Value = 1800
$ws.Cells.Item(8, 13).Value = 1800
$ws.Cells.Item(8, 16).Value = 1800

$ws.Cells.Item(9, 4).Value = 44503
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 10).Value = 72
$ws.Cells.Item(9, 11).Value = 1600
$ws.Cells.Item(9, 12).Value = 1600
$ws.Cells.Item(9, 13).Value = 1600
$ws.Cells.Item(9, 16).Value = 1600

$ws.Cells.Item(10, 4).Value = 44539
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 10).Value = 120
$ws.Cells.Item(10, 11).Value = 1600
$ws.Cells.Item(10, 12).Value = 1600
$ws.Cells.Item(10, 13).Value = 1600
$ws.Cells.Item(10, 16).Value = 1600

$ws.Cells.Item(11, 4).Value = 44176
$ws.Cells.Item(11, 8).Value = "Verde"
$ws.Cells.Item(11, 10).Value = 700
$ws.Cells.Item(11, 11).Value = 1600
$ws.Cells.Item(11, 12).Value = 1600
$ws.Cells.Item(11, 13).Value = 1600
$ws.Cells.Item(11, 16).Value = 1600

$ws.Cells.Item(12, 4).Value = 44159
$ws.Cells.Item(12, 8).Value = "Verde"
$ws.Cells.Item(12, 10).Value = 600
$ws.Cells.Item(12, 11).Value = 1600
$ws.Cells.Item(12, 12).Value = 1700
$ws.Cells.Item(12, 13).Value = 1650
$ws.Cells.Item(12, 16).Value = 1650

$ws.Cells.Item(13, 4).Value = 44494
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 10).Value = 120
$ws.Cells.Item(13, 11).Value = 1700
$ws.Cells.Item(13, 12).Value = 1700
$ws.Cells.Item(13, 13).Value = 1700
$ws.Cells.Item(13, 16).Value = 1700

$ws.Cells.Item(14, 4).Value = 44484
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 10).Value = 550
$ws.Cells.Item(14, 11).Value = 1700
$ws.Cells.Item(14, 12).Value = 1700
$ws.Cells.Item(14, 13).Value = 1700
$ws.Cells.Item(14, 16).Value = 1700

$ws.Cells.Item(15, 4).Value = 44488
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 10).Value = 600
$ws.Cells.Item(15, 11).Value = 1700
$ws.Cells.Item(15, 12).Value = 1800
$ws.Cells.Item(15, 13).Value = 1750
$ws.Cells.Item(15, 16).Value = 1750

$ws.Cells.Item(16, 4).Value = 44523
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 10).Value = 520
$ws.Cells.Item(16, 11).Value = 1800
$ws.Cells.Item(16, 12).Value = 1800
$ws.Cells.Item(16, 13).Value = 1800
$ws.Cells.Item(16, 16).Value = 1800

$ws.Cells.Item(17, 4).Value = 44540
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 10).Value = 600
$ws.Cells.Item(17, 11).Value = 1700
$ws.Cells.Item(17, 12).Value = 1700
$ws.Cells.Item(17, 13).Value = 1700
$ws.Cells.Item(17, 16).Value = 1700

$ws.Cells.Item(18, 4).Value = 44166
$ws.Cells.Item(18, 8).Value = "Verde"
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 1600
$ws.Cells.Item(18, 12).Value = 1600
$ws.Cells.Item(18, 13).Value = 1600
$ws.Cells.Item(18, 16).Value = 1600

$ws.Cells.Item(19, 4).Value = 44533
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 10).Value = 420
$ws.Cells.Item(19, 11).Value = 1700
$ws.Cells.Item(19, 12).Value = 1700
$ws.Cells.Item(19, 13).Value = 1700
$ws.Cells.Item(19, 16).Value = 1700

$ws.Cells.Item(20, 4).Value = 44481
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 10).Value = 300
$ws.Cells.Item(20, 11).Value = 1700
$ws.Cells.Item(20, 12).Value = 2000
$ws.Cells.Item(20, 13).Value = 1850
$ws.Cells.Item(20, 16).Value = 1850

$ws.Cells.Item(21, 4).Value = 44491
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 10).Value = 500
$ws.Cells.Item(21, 11).Value = 1700
$ws.Cells.Item(21, 12).Value = 1700
$ws.Cells.Item(21, 13).Value = 1700
$ws.Cells.Item(21, 16).Value = 1700

$ws.Cells.Item(22, 4).Value = 44511
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 10).Value = 144
$ws.Cells.Item(22, 11).Value = 1700
$ws.Cells.Item(22, 12).Value = 1700
$ws.Cells.Item(22, 13).Value = 1700
$ws.Cells.Item(22, 16).Value = 1700

$ws.Cells.Item(23, 4).Value = 44162
$ws.Cells.Item(23, 8).Value = "Verde"
$ws.Cells.Item(23, 10).Value = 700
$ws.Cells.Item(23, 11).Value = 1600
$ws.Cells.Item(23, 12).Value = 1600
$ws.Cells.Item(23, 13).Value = 1600
$ws.Cells.Item(23, 16).Value = 1600

$ws.Cells.Item(24, 4).Value = 44495
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 10).Value = 520
$ws.Cells.Item(24, 11).Value = 1800
$ws.Cells.Item(24, 12).Value = 1800
$ws.Cells.Item(24, 13).Value = 1800
$ws.Cells.Item(24, 16).Value = 1800

$ws.Cells.Item(25, 4).Value = 44161
$ws.Cells.Item(25, 8).Value = "Verde"
$ws.Cells.Item(25, 10).Value = 300
$ws.Cells.Item(25, 11).Value = 1700
$ws.Cells.Item(25, 12).Value = 1700
$ws.Cells.Item(25, 13).Value = 1700
$ws.Cells.Item(25, 16).Value = 1700

$ws.Cells.Item(26, 4).Value = 44498
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 10).Value = 500
$ws.Cells.Item(26, 11).Value = 1600
$ws.Cells.Item(26, 12).Value = 1600
$ws.Cells.Item(26, 13).Value = 1600
$ws.Cells.Item(26, 16).Value = 1600

$ws.Cells.Item(27, 4).Value = 44530
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 10).Value = 600
$ws.Cells.Item(27, 11).Value = 1500
$ws.Cells.Item(27, 12).Value = 1500
$ws.Cells.Item(27, 13).Value = 1500
$ws.Cells.Item(27, 16).Value = 1500

$ws.Cells.Item(28, 4).Value = 44519
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 10).Value = 600
$ws.Cells.Item(28, 11).Value = 1600
$ws.Cells.Item(28, 12).Value = 1800
$ws.Cells.Item(28, 13).Value = 1700
$ws.Cells.Item(28, 16).Value = 1700

$ws.Cells.Item(29, 4).Value = 44475
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 2000
$ws.Cells.Item(29, 12).Value = 2000
$ws.Cells.Item(29, 13).Value = 2000
$ws.Cells.Item(29, 16).Value = 2000

$ws.Cells.Item(30, 4).Value = 44516
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 10).Value = 360
$ws.Cells.Item(30, 11).Value = 1600
$ws.Cells.Item(30, 12).Value = 1600
$ws.Cells.Item(30, 13).Value = 1600
$ws.Cells.Item(30, 16).Value = 1600

$ws.Cells.Item(31, 4).Value = 44509
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 10).Value = 550
$ws.Cells.Item(31, 11).Value = 1700
$ws.Cells.Item(31, 12).Value = 1700
$ws.Cells.Item(31, 13).Value = 1700
$ws.Cells.Item(31, 16).Value = 1700

$ws.Cells.Item(32, 4).Value = 44168
$ws.Cells.Item(32, 8).Value = "Verde"
$ws.Cells.Item(32, 10).Value = 200
$ws.Cells.Item(32, 11).Value = 1600
$ws.Cells.Item(32, 12).Value = 1600
$ws.Cells.Item(32, 13).Value = 1600
$ws.Cells.Item(32, 16).Value = 1600

$ws.Cells.Item(33, 4).Value = 44490
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 10).Value = 72
$ws.Cells.Item(33, 11).Value = 1700
$ws.Cells.Item(33, 12).Value = 1700
$ws.Cells.Item(33, 13).Value = 1700
$ws.Cells.Item(33, 16).Value = 1700

$ws.Cells.Item(34, 4).Value = 44476
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 10).Value = 60
$ws.Cells.Item(34, 11).Value = 2000
$ws.Cells.Item(34, 12).Value = 2000
$ws.Cells.Item(34, 13).Value = 2000
$ws.Cells.Item(34, 16).Value = 2000

$ws.Cells.Item(35, 4).Value = 44529
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 10).Value = 120
$ws.Cells.Item(35, 11).Value = 1700
$ws.Cells.Item(35, 12).Value = 1700
$ws.Cells.Item(35, 13).Value = 1700
$ws.Cells.Item(35, 16).Value = 1700

$ws.Cells.Item(36, 4).Value = 44517
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 10).Value = 120
$ws.Cells.Item(36, 11).Value = 1600
$ws.Cells.Item(36, 12).Value = 1600
$ws.Cells.Item(36, 13).Value = 1600
$ws.Cells.Item(36, 16).Value = 1600

$ws.Cells.Item(37, 4).Value = 44526
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 10).Value = 600
$ws.Cells.Item(37, 11).Value = 1700
$ws.Cells.Item(37, 12).Value = 1700
$ws.Cells.Item(37, 13).Value = 1700
$ws.Cells.Item(37, 16).Value = 1700

$ws.Cells.Item(38, 4).Value = 44179
$ws.Cells.Item(38, 8).Value = "Verde"
$ws.Cells.Item(38, 10).Value = 200
$ws.Cells.Item(38, 11).Value = 1600
$ws.Cells.Item(38, 12).Value = 1600
$ws.Cells.Item(38, 13).Value = 1600
$ws.Cells.Item(38, 16).Value = 1600

$ws.Cells.Item(39, 4).Value = 44473
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 10).Value = 100
$ws.Cells.Item(39, 11).Value = 2000
$ws.Cells.Item(39, 12).Value = 2000
$ws.Cells.Item(39, 13).Value = 2000
$ws.Cells.Item(39, 16).Value = 2000

$ws.Cells.Item(40, 4).Value = 44544
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 10).Value = 300
$ws.Cells.Item(40, 11).Value = 1700
$ws.Cells.Item(40, 12).Value = 1700
$ws.Cells.Item(40, 13).Value = 1700
$ws.Cells.Item(40, 16).Value = 1700

$ws.Cells.Item(41, 4).Value = 44537
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 10).Value = 480
$ws.Cells.Item(41, 11).Value = 1600
$ws.Cells.Item(41, 12).Value = 1600
$ws.Cells.Item(41, 13).Value = 1600
$ws.Cells.Item(41, 16).Value = 1600

$ws.Cells.Item(42, 4).Value = 44482
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 10).Value = 72
$ws.Cells.Item(42, 11).Value = 2000
$ws.Cells.Item(42, 12).Value = 2000
$ws.Cells.Item(42, 13).Value = 2000
$ws.Cells.Item(42, 16).Value = 2000

$ws.Cells.Item(43, 4).Value = 44165
$ws.Cells.Item(43, 8).Value = "Verde"
$ws.Cells.Item(43, 10).Value = 300
$ws.Cells.Item(43, 11).Value = 1600
$ws.Cells.Item(43, 12).Value = 1600
$ws.Cells.Item(43, 13).Value = 1600
$ws.Cells.Item(43, 16).Value = 1600

$ws.Cells.Item(44, 4).Value = 44496
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 10).Value = 84
$ws.Cells.Item(44, 11).Value = 1800
$ws.Cells.Item(44, 12).Value = 1800
$ws.Cells.Item(44, 13).Value = 1800
$ws.Cells.Item(44, 16).Value = 1800

$ws.Cells.Item(45, 4).Value = 44504
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 10).Value = 180
$ws.Cells.Item(45, 11).Value = 1600
$ws.Cells.Item(45, 12).Value = 1600
$ws.Cells.Item(45, 13).Value = 1600
$ws.Cells.Item(45, 16).Value = 1600

$ws.Cells.Item(46, 4).Value = 44522
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 10).Value = 180
$ws.Cells.Item(46, 11).Value = 1800
$ws.Cells.Item(46, 12).Value = 1800
$ws.Cells.Item(46, 13).Value = 1800
$ws.Cells.Item(46, 16).Value = 1800
